# Update the "Donnees COVID-19 Valais" daily tracking sheet.
#
# Only the literal input columns (C, E, F, G, L, M) are touched; columns
# B, H, J, K are "TODAY()"-driven shared formulas (si=43..46) that
# recalculate automatically from those inputs once the workbook recalcs.
#
# Columns L ("Nb nouveaux deces a l'hopital") and M ("Nb nouveaux deces
# extra-hospitaliers") carry a Text ("@") number format even though they
# only ever hold small integers. Writing a plain number into a Text-
# formatted cell via Range.Value (or PasteSpecial -xlPasteValues, or a
# bulk array write) gets stored as a text string, which does not match
# the target file (plain numeric <v> cells). To write a genuine number
# while preserving the existing cell style, we temporarily borrow the
# General-formatted style from the neighbouring D/E column (same border,
# General number format), assign the numeric value, then restore the
# original L/M style by copying formats back from an untouched reference
# cell further up the same column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumericTextCell($row, $col, $generalSourceRow, $generalSourceCol, $styleSourceRow, $styleSourceCol, $value) {
    $ws.Cells.Item($generalSourceRow, $generalSourceCol).Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Cells.Item($row, $col).Value = $value
    $ws.Cells.Item($styleSourceRow, $styleSourceCol).Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4122)   # xlPasteFormats
}

# --- Row 594: new-case count revised upward, plus 1 new extra-hospital death ---
$ws.Cells.Item(594, 3).Value = 61   # C594 Nb nouveaux cas positifs
Set-NumericTextCell 594 12 594 4 566 12 1   # L594 Nb nouveaux deces a l'hopital: 0 -> 1

# --- Row 602: new-case count revised ---
$ws.Cells.Item(602, 3).Value = 27   # C602

# --- Row 603: new-case count revised ---
$ws.Cells.Item(603, 3).Value = 60   # C603

# --- Row 604: new-case count + hospitalisation count revised ---
$ws.Cells.Item(604, 3).Value = 37   # C604
$ws.Cells.Item(604, 7).Value = 9    # G604 Total hospitalisations COVID-19

# --- Row 605: previously-blank day, now filled in with real data ---
$ws.Cells.Item(605, 3).Value = 43   # C605
$ws.Cells.Item(605, 5).Value = 3    # E605
$ws.Cells.Item(605, 6).Value = 2    # F605
$ws.Cells.Item(605, 7).Value = 8    # G605
Set-NumericTextCell 605 12 605 4 566 12 0   # L605
Set-NumericTextCell 605 13 605 5 566 13 0   # M605

# --- Row 606: previously-blank day, now filled in with real data ---
$ws.Cells.Item(606, 3).Value = 27   # C606
$ws.Cells.Item(606, 5).Value = 2    # E606
$ws.Cells.Item(606, 6).Value = 2    # F606
$ws.Cells.Item(606, 7).Value = 9    # G606
Set-NumericTextCell 606 12 606 4 566 12 0   # L606
Set-NumericTextCell 606 13 606 5 566 13 0   # M606

# --- Row 607: previously-blank day, now filled in with real data ---
$ws.Cells.Item(607, 3).Value = 14   # C607
$ws.Cells.Item(607, 5).Value = 2    # E607
$ws.Cells.Item(607, 6).Value = 2    # F607
$ws.Cells.Item(607, 7).Value = 8    # G607
Set-NumericTextCell 607 12 607 4 566 12 0   # L607
Set-NumericTextCell 607 13 607 5 566 13 0   # M607

# --- Row 608: previously-blank day, now filled in with real data ---
$ws.Cells.Item(608, 3).Value = 3    # C608
$ws.Cells.Item(608, 5).Value = 2    # E608
$ws.Cells.Item(608, 6).Value = 2    # F608
$ws.Cells.Item(608, 7).Value = 8    # G608
Set-NumericTextCell 608 12 608 4 566 12 0   # L608
Set-NumericTextCell 608 13 608 5 566 13 0   # M608

$wb.Application.CutCopyMode = $false
